# Rapport du 22 Septembre 2025
# Applies the workbook update:
#  - widen column M (13) to match column L's width (new <col> override)
#  - correct row 45's quantity/amount (Cafe Altimo 50g carton count)
#  - add "Stationnement" (parking) cost of 6000 for every "Vente" row
#    between rows 54 and 94, and bump row 85's quantity/amount as well

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column width: give column M (Stationnement) the same display width
#     as column L (Transport) already has. The ColumnWidth COM property is
#     expressed in "characters", which Excel internally offsets by 5/6 of a
#     character versus the width stored in the sheet XML, so we subtract
#     that offset before writing it.
$ws.Columns.Item(13).ColumnWidth = 6.01171875 - 0.8333333333333334

# --- row 45: revised carton quantity & amount for the Cafe Altimo 50g line
$ws.Range("J45").Value = 26.291666
$ws.Range("K45").Value = 460104.16

# --- rows 54-94: add the 6000 "Stationnement" charge that was missing
for ($r = 54; $r -le 94; $r++) {
    $ws.Range("M" + $r).Value = 6000
}

# --- row 85: quantity & amount were also corrected alongside the new charge
$ws.Range("J85").Value = 0.7
$ws.Range("K85").Value = 8575.0
